$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Layer0")
$ws1.Range("B2").Value = -0.4400819820205796
$ws1.Range("C2").Value = -0.1827760305871028
$ws1.Range("B3").Value = -1.103531195188436
$ws1.Range("C3").Value = -0.2766336750638603
$ws1.Range("B4").Value = 0.3074702297549578
$ws1.Range("C4").Value = -1.53290043198063

$ws2 = $wb.Worksheets.Item("Layer1")
$ws2.Range("B2").Value = -0.4356415220920893
$ws2.Range("C2").Value = -0.3284568562298597
$ws2.Range("B3").Value = -0.5743275939144695
$ws2.Range("C3").Value = 0.6744243603619372
$ws2.Range("B4").Value = -1.908292008488429
$ws2.Range("C4").Value = 0.4097478028414245
